# Fix missing value in test ARC Assay file
# The "Proteomix" sheet is missing the "Parameter [protein modification]"
# value ("N14-oxyginated") for the third data row (row 4, sample3), while
# rows 2 and 3 already have it filled in. Also update the active cell
# selection to reflect the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proteomix")

$ws.Range("E4").Value = "N14-oxyginated"

$ws.Activate()
$ws.Range("E4").Select()
